# ---- Step 1: Insert new "2022-Q4" sheet before "2022-Q3" ----
$wb = $excel.ActiveWorkbook
$templateSheet = $wb.Worksheets.Item(2)   # "2022-Q3" sheet (becomes index 3 after insert)
$ws = $wb.Worksheets.Add($templateSheet)   # inserted BEFORE $templateSheet -> new index 2
$ws.Name = "2022-Q4"

# ---- Step 2: header row (row 1), columns B..H ----
$ws.Range("B1:H1").Font.Bold = $true
$ws.Range("B1:H1").HorizontalAlignment = -4108
$ws.Range("B1:H1").VerticalAlignment = -4160
$ws.Range("B1:H1").Borders.LineStyle = 1
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# ---- Step 3: data rows 2..20 ----
$ws.Range("A2:A20").Font.Bold = $true
$ws.Range("A2:A20").HorizontalAlignment = -4108
$ws.Range("A2:A20").VerticalAlignment = -4160
$ws.Range("A2:A20").Borders.LineStyle = 1
$ws.Range("B2:G20").NumberFormat = "@"   # keep fund codes/figures as text (preserve leading zeros)

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "010379"
$ws.Cells.Item(2,3).Value = "广发均衡优选混合A"
$ws.Cells.Item(2,4).Value = "48.69"
$ws.Cells.Item(2,5).Value = "64.94"
$ws.Cells.Item(2,6).Value = "4.51"
$ws.Cells.Item(2,7).Value = "2.1959"
$ws.Cells.Item(2,8).Value = 9
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "009887"
$ws.Cells.Item(3,3).Value = "广发稳健优选六个月持有期混合A"
$ws.Cells.Item(3,4).Value = "17.19"
$ws.Cells.Item(3,5).Value = "64.90"
$ws.Cells.Item(3,6).Value = "5.00"
$ws.Cells.Item(3,7).Value = "0.8595"
$ws.Cells.Item(3,8).Value = 9
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "270022"
$ws.Cells.Item(4,3).Value = "广发内需增长混合A"
$ws.Cells.Item(4,4).Value = "9.81"
$ws.Cells.Item(4,5).Value = "79.97"
$ws.Cells.Item(4,6).Value = "6.30"
$ws.Cells.Item(4,7).Value = "0.6180"
$ws.Cells.Item(4,8).Value = 7
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "009888"
$ws.Cells.Item(5,3).Value = "广发稳健优选六个月持有期混合C"
$ws.Cells.Item(5,4).Value = "10.60"
$ws.Cells.Item(5,5).Value = "64.90"
$ws.Cells.Item(5,6).Value = "5.00"
$ws.Cells.Item(5,7).Value = "0.5300"
$ws.Cells.Item(5,8).Value = 9
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "011134"
$ws.Cells.Item(6,3).Value = "广发价值优选混合A"
$ws.Cells.Item(6,4).Value = "4.58"
$ws.Cells.Item(6,5).Value = "94.24"
$ws.Cells.Item(6,6).Value = "7.70"
$ws.Cells.Item(6,7).Value = "0.3527"
$ws.Cells.Item(6,8).Value = 5
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "010380"
$ws.Cells.Item(7,3).Value = "广发均衡优选混合C"
$ws.Cells.Item(7,4).Value = "3.21"
$ws.Cells.Item(7,5).Value = "64.94"
$ws.Cells.Item(7,6).Value = "4.51"
$ws.Cells.Item(7,7).Value = "0.1448"
$ws.Cells.Item(7,8).Value = 9
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "011183"
$ws.Cells.Item(8,3).Value = "广发内需增长混合C"
$ws.Cells.Item(8,4).Value = "1.82"
$ws.Cells.Item(8,5).Value = "79.97"
$ws.Cells.Item(8,6).Value = "6.30"
$ws.Cells.Item(8,7).Value = "0.1147"
$ws.Cells.Item(8,8).Value = 7
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "011135"
$ws.Cells.Item(9,3).Value = "广发价值优选混合C"
$ws.Cells.Item(9,4).Value = "1.38"
$ws.Cells.Item(9,5).Value = "94.24"
$ws.Cells.Item(9,6).Value = "7.70"
$ws.Cells.Item(9,7).Value = "0.1063"
$ws.Cells.Item(9,8).Value = 5
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "002295"
$ws.Cells.Item(10,3).Value = "广发稳安灵活配置混合A"
$ws.Cells.Item(10,4).Value = "1.51"
$ws.Cells.Item(10,5).Value = "80.18"
$ws.Cells.Item(10,6).Value = "4.60"
$ws.Cells.Item(10,7).Value = "0.0695"
$ws.Cells.Item(10,8).Value = 6
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "003401"
$ws.Cells.Item(11,3).Value = "工银可转债债券"
$ws.Cells.Item(11,4).Value = "3.21"
$ws.Cells.Item(11,5).Value = "49.19"
$ws.Cells.Item(11,6).Value = "2.07"
$ws.Cells.Item(11,7).Value = "0.0664"
$ws.Cells.Item(11,8).Value = 8
$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "002076"
$ws.Cells.Item(12,3).Value = "浙商中证500指数增强A"
$ws.Cells.Item(12,4).Value = "6.56"
$ws.Cells.Item(12,5).Value = "87.04"
$ws.Cells.Item(12,6).Value = "1.01"
$ws.Cells.Item(12,7).Value = "0.0663"
$ws.Cells.Item(12,8).Value = 8
$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "005945"
$ws.Cells.Item(13,3).Value = "工银可转债优选债券A"
$ws.Cells.Item(13,4).Value = "2.67"
$ws.Cells.Item(13,5).Value = "38.35"
$ws.Cells.Item(13,6).Value = "1.82"
$ws.Cells.Item(13,7).Value = "0.0486"
$ws.Cells.Item(13,8).Value = 8
$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "009956"
$ws.Cells.Item(14,3).Value = "广发恒誉混合A"
$ws.Cells.Item(14,4).Value = "2.44"
$ws.Cells.Item(14,5).Value = "25.40"
$ws.Cells.Item(14,6).Value = "1.40"
$ws.Cells.Item(14,7).Value = "0.0342"
$ws.Cells.Item(14,8).Value = 9
$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "009649"
$ws.Cells.Item(15,3).Value = "嘉实精选平衡混合A"
$ws.Cells.Item(15,4).Value = "0.47"
$ws.Cells.Item(15,5).Value = "68.05"
$ws.Cells.Item(15,6).Value = "6.21"
$ws.Cells.Item(15,7).Value = "0.0292"
$ws.Cells.Item(15,8).Value = 4
$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "007386"
$ws.Cells.Item(16,3).Value = "浙商中证500指数增强C"
$ws.Cells.Item(16,4).Value = "1.70"
$ws.Cells.Item(16,5).Value = "87.04"
$ws.Cells.Item(16,6).Value = "1.01"
$ws.Cells.Item(16,7).Value = "0.0172"
$ws.Cells.Item(16,8).Value = 8
$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "005946"
$ws.Cells.Item(17,3).Value = "工银可转债优选债券C"
$ws.Cells.Item(17,4).Value = "0.75"
$ws.Cells.Item(17,5).Value = "38.35"
$ws.Cells.Item(17,6).Value = "1.82"
$ws.Cells.Item(17,7).Value = "0.0136"
$ws.Cells.Item(17,8).Value = 8
$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "009650"
$ws.Cells.Item(18,3).Value = "嘉实精选平衡混合C"
$ws.Cells.Item(18,4).Value = "0.04"
$ws.Cells.Item(18,5).Value = "68.05"
$ws.Cells.Item(18,6).Value = "6.21"
$ws.Cells.Item(18,7).Value = "0.0025"
$ws.Cells.Item(18,8).Value = 4
$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "009957"
$ws.Cells.Item(19,3).Value = "广发恒誉混合C"
$ws.Cells.Item(19,4).Value = "0.06"
$ws.Cells.Item(19,5).Value = "25.40"
$ws.Cells.Item(19,6).Value = "1.40"
$ws.Cells.Item(19,7).Value = "0.0008"
$ws.Cells.Item(19,8).Value = 9
$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "008604"
$ws.Cells.Item(20,3).Value = "广发稳安灵活配置混合C"
$ws.Cells.Item(20,4).Value = "0.01"
$ws.Cells.Item(20,5).Value = "80.18"
$ws.Cells.Item(20,6).Value = "4.60"
$ws.Cells.Item(20,7).Value = "0.0005"
$ws.Cells.Item(20,8).Value = 6

# ---- Step 4: update the "总计" (summary) sheet: insert a 2022-Q4 row at row 2 ----
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Cells.Item(3,1).Copy($summary.Cells.Item(2,1))   # copy A-column index style onto the new row
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 19
$summary.Cells.Item(2,4).Value = 5.27
$summary.Range("B2:D2").ClearFormats()

# re-number the "A" index column (0-based rank) for the rows pushed down by the insert
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(7,1).Value = 5
